# "Correction engine motor + overflow bottles"
#
# - Add the two missing "Description" notes for the engine assembly rows
#   (Honda CBR600RR Engine / Thermostat).
# - Resize the Description column (E) a bit narrower and the Make/Buy
#   column (D) a bit wider to fit the new / updated text, and let the
#   now-shorter overflow-bottle descriptions (rows 18-19) autofit back to
#   the default row height instead of staying tall.
# - Grow rows 15-16 (shifter axis / shifter gear) so their long
#   descriptions are fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New engine notes
$ws.Range("E3").Value = "Second-hand, PC40"
$ws.Range("E7").Value = "From PC37, outlet at the opposite of PC40 version"

# Column widths: D a little wider, E a little narrower (no longer "best fit")
$ws.Columns(4).ColumnWidth = 23.5
$ws.Columns(5).ColumnWidth = 53

# Taller rows so the shifter axis/gear descriptions aren't clipped
$ws.Rows(15).RowHeight = 29.4
$ws.Rows(16).RowHeight = 30

# Overflow bottle rows go back to the default (auto) row height
$ws.Rows(18).AutoFit() | Out-Null
$ws.Rows(19).AutoFit() | Out-Null

# Leave the cursor where the edit was made
$ws.Range("E10").Select() | Out-Null
